# Actualización automática 2025-09-30 17:30:09
# Updates sales figures for "RIOS CARRION ANGEL BENIGNO" across the three
# sheets of the workbook (per-group sales, monthly sales, monthly
# compliance), reflecting newly recorded September sales for clients
# "CONSTANTE CAMACHO ARIANA ELIZABETH" (INODOROS + PORCELANATO) and
# "CULMA OVIEDO NINI JOHANA" (240X80 PORCELANATO).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("H8").Value = 1128.6
$wsGrupo.Range("M8").Value = 1283.09
$wsGrupo.Range("D10").Value = 549.5

$wsGrupo.Range("D26").Value = "1 de 24"
$wsGrupo.Range("H26").Value = "2 de 24"
$wsGrupo.Range("M26").Value = "7 de 24"

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F8").Value = 2411.69
$wsMensual.Range("F10").Value = 549.5
$wsMensual.Range("F26").Value = 23879.02

# ---------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column widths: the stored OOXML width is ColumnWidth + 5/6, so back
# out the offset to land on the exact target widths (13 and 24).
$wsCumpl.Columns.Item(4).ColumnWidth = 13 - (5 / 6)
$wsCumpl.Columns.Item(5).ColumnWidth = 24 - (5 / 6)

# Row 3 - 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 549.5
$wsCumpl.Range("E3").Value = 6925.9083879616
$wsCumpl.Range("F3").Value = 0.07350768967818734

# Row 6 - INODOROS
$wsCumpl.Range("D6").Value = 1434.84
$wsCumpl.Range("E6").Value = -527.6738913843989
$wsCumpl.Range("F6").Value = 1.581672845108451

# Row 12 - PORCELANATO
$wsCumpl.Range("D12").Value = 21949
$wsCumpl.Range("E12").Value = 21151.0854117774
$wsCumpl.Range("F12").Value = 0.5092565314035847

# Row 15 - TOTAL
$wsCumpl.Range("D15").Value = 23879.02
$wsCumpl.Range("E15").Value = 34324.44623249459
$wsCumpl.Range("F15").Value = 0.410268005424538
